$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for A2:A37 (replacing the previous random-like series)
$values = @(
    0.60995931606694376,
    0.83355497761059105,
    0.73250076650668305,
    0.7165717039853946,
    0.61984130807822369,
    0.66931913000456766,
    0.72186787213185144,
    0.65183093080878007,
    0.64675694385725657,
    0.66252816554234628,
    0.69000103550874581,
    0.66309721044387027,
    0.62459718320523083,
    0.64702839335258366,
    0.6550904217207415,
    0.67643270861422555,
    0.68753693218394185,
    0.67552390375129034,
    0.67240435903540374,
    0.62937618625795244,
    0.64890432386180164,
    0.59441390974804431,
    0.61656945535666285,
    0.70441661544026879,
    0.65352193853913143,
    0.62280237368780889,
    0.72343567692511934,
    0.70057597564136309,
    0.66098272492859622,
    0.63271021164838903,
    0.65393952664769461,
    0.79501311532159047,
    0.74672234775261737,
    0.63240445833042613,
    0.7167552527227784,
    2.1071601324810114
)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 1).Value = $v
    $row++
}

# Narrow column A slightly (from 13.7109375 down to 12.7109375 characters)
$ws.Range("A1:A1").EntireColumn.ColumnWidth = 11.8
